$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value = "Datos actualizados a 22 de Octubre de 2020 a las 00:09"

# Row 4
$ws.Range("B4").Value = 8571906
$ws.Range("C4").Value = 50956
$ws.Range("D4").Value = 5575461
$ws.Range("E4").Value = 2769304
$ws.Range("G4").Value = 957
$ws.Range("H4").Value = 227141

# Row 6
$ws.Range("B6").Value = 5298772
$ws.Range("C6").Value = 23955
$ws.Range("E6").Value = 421777
$ws.Range("G6").Value = 514
$ws.Range("H6").Value = 155402

# Row 10
$ws.Range("B10").Value = 981700
$ws.Range("C10").Value = 7561
$ws.Range("D10").Value = 884895
$ws.Range("E10").Value = 67341
$ws.Range("G10").Value = 192
$ws.Range("H10").Value = 29464

# Row 48
$ws.Range("B48").Value = 105883
$ws.Range("C48").Value = 178
$ws.Range("D48").Value = 98516
$ws.Range("E48").Value = 1212
$ws.Range("G48").Value = 13
$ws.Range("H48").Value = 6155

# Row 49
$ws.Range("B49").Value = 102415
$ws.Range("C49").Value = 196
$ws.Range("D49").Value = 92149
$ws.Range("E49").Value = 6699
$ws.Range("G49").Value = 21
$ws.Range("H49").Value = 3567

# Row 84
$ws.Range("A84").Value = "Bulgaria"
$ws.Range("B84").Value = 33335
$ws.Range("C84").Value = 1472
$ws.Range("D84").Value = 17598
$ws.Range("E84").Value = 14689
$ws.Range("G84").Value = 29
$ws.Range("H84").Value = 1048

# Row 85
$ws.Range("A85").Value = "El Salvador"
$ws.Range("B85").Value = 32120
$ws.Range("C85").Value = 145
$ws.Range("D85").Value = 27670
$ws.Range("E85").Value = 3517
$ws.Range("G85").Value = 4
$ws.Range("H85").Value = 933

# Row 120
$ws.Range("B120").Value = 7638
$ws.Range("C120").Value = 4
$ws.Range("D120").Value = 7363

# Row 146
$ws.Range("A146").Value = "Guyana"
$ws.Range("B146").Value = 3850
$ws.Range("C146").Value = 54
$ws.Range("D146").Value = 2839
$ws.Range("E146").Value = 895
$ws.Range("G146").Value = 2
$ws.Range("H146").Value = 116

# Row 147
$ws.Range("A147").Value = "Principado de Andorra"
$ws.Range("B147").Value = 3811
$ws.Range("D147").Value = 2470
$ws.Range("E147").Value = 1278
$ws.Range("G147").Value = 1
$ws.Range("H147").Value = 63

# Row 148
$ws.Range("A148").Value = "Letonia"
$ws.Range("B148").Value = 3797
$ws.Range("C148").Value = 188
$ws.Range("D148").Value = 1341
$ws.Range("E148").Value = 2409
$ws.Range("H148").Value = 47

# Row 161
$ws.Range("B161").Value = 2120
$ws.Range("C161").Value = 16
$ws.Range("D161").Value = 1561
$ws.Range("E161").Value = 508

# Row 167
$ws.Range("B167").Value = 1214
$ws.Range("C167").Value = 2
$ws.Range("E167").Value = 17

# Row 172
$ws.Range("B172").Value = 762
$ws.Range("C172").Value = 6
$ws.Range("E172").Value = 65
